$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.424.99"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.850.87"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.14"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4758"
$ws.Range("E7").Value = "  +3.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2747"
$ws.Range("E8").Value = "  +2.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06321"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.56"
$ws.Range("E10").Value = "  +9.37%  "

$ws.Range("D11").Value = "1.874.40"
$ws.Range("E11").Value = "  +2.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07458"
$ws.Range("E12").Value = "  +1.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.947"
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.69"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6242"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").Value = "30.387.06"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "246.35"
$ws.Range("E17").Value = "  +8.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  +2.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007320"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.902"
$ws.Range("E22").Value = "  +1.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.903"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.59"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.083"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.97"
$ws.Range("E26").Value = "  +1.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.869"
$ws.Range("E27").Value = "  +1.75%  "

$ws.Range("E28").Value = "  +1.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.347"
$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.031"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.816"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04839"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6959"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.702"
$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01900"
$ws.Range("E36").Value = "  +5.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  +3.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8773"
$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.995"
$ws.Range("E39").Value = "  +4.72%  "

$ws.Range("E40").Value = "  +4.13%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4050"
$ws.Range("E42").Value = "  +1.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.499"
$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.171"
$ws.Range("E44").Value = "  +4.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.03"
$ws.Range("E45").Value = "  +6.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1197"
$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.68"
$ws.Range("E47").Value = "  +3.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.523"
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05503"
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.349"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3678"
$ws.Range("E51").Value = "  +1.76%  "
